# Update column G ("K") values on the active sheet to reflect the
# regenerated save_data (K computed instead of Strike#, std/mean
# recalculated, s_vals written).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$gValues = @{
    2  = 3
    3  = 2
    4  = 1
    5  = 0
    6  = 2
    7  = 1
    8  = 1
    9  = 1
    10 = 1
    11 = 0
    12 = 0
    13 = 0
    14 = 0
    15 = 0
    16 = 3
    17 = 1
    18 = 1
    19 = 3
    20 = 2
    21 = 1
    22 = 1
    23 = 2
    24 = 1
    25 = 0
    26 = 1
    27 = 1
    28 = 0
    29 = 0
    30 = 2
    31 = 2
    32 = 1
    33 = 1
    34 = 4
    35 = 3
    36 = 1
    37 = 2
    38 = 3
    39 = 3
    40 = 1
}

foreach ($row in $gValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $gValues[$row]
}
